$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 describes the "Encode nominal to numerical" feature. Rename it to
# "Create indicator variables", mark it complete, and fill in the
# "Value to user" column with the same text as the feature name.
$ws.Range("A9").Value = "Create indicator variables"
$ws.Range("D9").Value = "COMPLETE"
$ws.Range("E9").Value = "Create indicator variables"

# Move the active selection from J8 to A9.
$ws.Range("A9").Select()
